$wb = $excel.ActiveWorkbook

# Work on the "Repayment schedule" sheet: insert a new column before column N
# (shifting the old "Late" / "heading" / "Outstanding" columns one to the right).
$ws = $wb.Worksheets.Item("Repayment schedule")

$ws.Columns("N:N").Insert()

# Match the width Excel assigns to a freshly inserted column (copies the
# width of the column to its left, "In Advance" = stored width 11).
$ws.Columns("N:N").ColumnWidth = 10.166666666666666

# Make "Repayment schedule" the active sheet/tab and select cell R8 on it,
# which also clears the previous tab-selection on "Edit Repayment Schedule".
$ws.Activate()
$ws.Range("R8").Select()
